$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# First, touch one cell per *new* unique string, in the exact order the
# strings need to be appended to the shared-strings table so that the
# resulting xl/sharedStrings.xml ordering (and therefore numeric <v> indexes
# used throughout the sheet) matches the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("G45").Value = "Filler"                          # new string 49
$ws.Range("H45").Value = "Proposed"                         # new string 50
$ws.Range("F48").Value = "Accepted"                          # new string 51
$ws.Range("B46").Value = "< Filler"                          # new string 52
$ws.Range("B50").Value = "> Filler"                          # new string 53
$ws.Range("B55").Value = "'= Filler"                         # new string 54 (leading '=' needs quote-prefix)
$ws.Range("B60").Value = "No Filler"                          # new string 55
$ws.Range("C47").Value = "And pie has filler"                 # new string 56
$ws.Range("C61").Value = "And pie has NO filler"               # new string 57
$ws.Range("D60").Value = "an ingrediuient % is changed"        # new string 58
$ws.Range("E60").Value = "should reject the proposal"          # new string 59

# ---------------------------------------------------------------------------
# Now fill in the rest of the new block (rows 45-61). Cells referencing
# strings already used above will simply reuse those shared-string indexes.
# ---------------------------------------------------------------------------

# --- Row 46 ---
$ws.Range("C46").Value = "Ingredients are added"
$ws.Range("F46").Value = "Before"
$ws.Range("G46").Value = 40
$ws.Range("H46").Value = 20

# --- Row 47 ---
$ws.Range("F47").Value = "After"
$ws.Range("G47").Value = 20

# --- Row 48 ---
$ws.Range("H48").Value = 20

# --- Row 50 ---
$ws.Range("C50").Value = "Ingredients are added"
$ws.Range("G50").Value = "Filler"
$ws.Range("H50").Value = "Proposed"

# --- Row 51 ---
$ws.Range("C51").Value = "And pie has filler"
$ws.Range("F51").Value = "Before"
$ws.Range("G51").Value = 40
$ws.Range("H51").Value = 60

# --- Row 52 ---
$ws.Range("F52").Value = "After"
$ws.Range("G52").Value = 0

# --- Row 53 ---
$ws.Range("F53").Value = "Accepted"
$ws.Range("H53").Value = 20

# --- Row 55 ---
$ws.Range("C55").Value = "Ingredients are added"
$ws.Range("G55").Value = "Filler"
$ws.Range("H55").Value = "Proposed"

# --- Row 56 ---
$ws.Range("C56").Value = "And pie has filler"
$ws.Range("F56").Value = "Before"
$ws.Range("G56").Value = 40
$ws.Range("H56").Value = 40

# --- Row 57 ---
$ws.Range("F57").Value = "After"
$ws.Range("G57").Value = 0

# --- Row 58 ---
$ws.Range("F58").Value = "Accepted"
$ws.Range("H58").Value = 40

# --- Row 60 ---
$ws.Range("C60").Value = "Ingredients are added"

# --- View changes ---
# (topLeftCell scroll position is window chrome state that this host does not
#  persist back to the saved sheetView; Zoom/selection are still applied.)
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 2
$win.Zoom = 100
$ws.Range("E60").Select() | Out-Null
